# task9&10_v2.0.2.xlsx update
# "diary updated, sprint backlog, alertpanel angepasst"
#
# Product Backlog: stories 1-3 are now finished; fill in the actual effort
# (Effort Actual) for each and flip Status from "not started" to "finished".
#
# Sprint Backlog: the remaining "pending" stories are now "finished" too.

$wb = $excel.ActiveWorkbook

# --- Product Backlog -------------------------------------------------
$productBacklog = $wb.Worksheets.Item("Product Backlog")

# Story 1 (release alert)
$productBacklog.Range("G2").Value = "47h"
$productBacklog.Range("H2").Value = "finished"

# Story 2 (patient referral)
$productBacklog.Range("G3").Value = "70h"
$productBacklog.Range("H3").Value = "finished"

# Story 3 (denial handling)
$productBacklog.Range("F4").Value = "45h"
$productBacklog.Range("G4").Value = "60h"
$productBacklog.Range("H4").Value = "finished"

# --- Sprint Backlog ----------------------------------------------------
$sprintBacklog = $wb.Worksheets.Item("Sprint Backlog")

$sprintBacklog.Range("K8").Value = "finished"
$sprintBacklog.Range("K9").Value = "finished"
$sprintBacklog.Range("K18").Value = "finished"
$sprintBacklog.Range("K19").Value = "finished"
$sprintBacklog.Range("K26").Value = "finished"
$sprintBacklog.Range("K27").Value = "finished"

# --- Selection / active sheet, matching the author's final view -------
$sprintBacklog.Range("I25").Select()

$productBacklog.Activate()
$productBacklog.Range("G5").Select()
